# Regenerate the handoff report: a new GUID-named file was produced (replacing the
# previous d36a25b3-... one), new content hashes for the generated xliff files, and
# fresh handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "7d0caaa4-784d-4383-8fcd-afeb58712726"
$newHash = "a8b938be2a9090e4993f6e7ff7df2b3297997b78"

$newFileName   = "$newGuid.md"
$newDisplay    = "e2e\" + $newFileName

$newHoDate     = "2016-09-06 03:05:17"
$newZhXlfDate  = "2016-09-06 03:05:13"

$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

# The external hyperlink target (github blob URL) is unchanged - only the
# human-readable text shown for the link moves to the new file name.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba8c39202028c60395f7af1ff1d5cc5786f97889/e2e/d36a25b3-6916-4364-a7ea-d72c0fce0277.md"

function Update-Hyperlink($ws, $cellRef, $display) {
    $range = $ws.Range($cellRef)
    $range.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($range, $linkAddress, $null, $null, $display)
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newDisplay
Update-Hyperlink $wsOverview "B2" $newDisplay
$wsOverview.Range("G2").Value = $newHoDate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
Update-Hyperlink $wsZh "A2" $newFileName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhXlfDate

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
Update-Hyperlink $wsDe "A2" $newFileName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoDate
